$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175886750221252
$ws.Range("B1").Value = 2.138853549957275
$ws.Range("C1").Value = 3.294553518295288
$ws.Range("D1").Value = 2.226548433303833
$ws.Range("E1").Value = 0.7578895092010498
